# 6th testing of OB using 1419 data, 10 epoch, and do the spliting data to data train & test
$wb = $excel.ActiveWorkbook

# The "Real" sheet is the first sheet in the workbook (tab name "Real").
$ws = $wb.Worksheets.Item(1)

# Fill in the confusion-matrix counts for the new 6th test run (row 8).
$ws.Range("H8").Value = 123
$ws.Range("I8").Value = 8
$ws.Range("J8").Value = 22
$ws.Range("K8").Value = 0

# Record the saved training/testing run folders for this test.
$ws.Range("O8").Value = "runs\detect\train23"
$ws.Range("P8").Value = "runs\detect\train232"

# Leave the selection where the author last left it when saving.
$ws.Range("H12").Select() | Out-Null
